# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" (column G) cells so that entries read
# "System, dnasr281@gmail.com" instead of "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
